# OLX Monitor 2026-02-17 17:05
# Appends 7 newly-discovered listings to the bottom of the "PODSUMOWANIE"
# sheet's per-listing table (rows 54-60), extending its used range from
# A1:H53 to A1:H60.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PODSUMOWANIE")

# Copy the formatting of the last existing data row (53) down onto the new
# rows so they inherit the same cell styles used by the rest of the table.
$ws.Range("A53:H53").Copy()
$ws.Range("A54:H60").PasteSpecial(-4122)

# Rows 54 and 58 are "days on market = 29" listings, which this sheet
# highlights with the lighter "F" style (same one used on row 47/51) rather
# than the "long-listed" style the rest of column F got from row 53 above.
$ws.Range("F47").Copy()
$ws.Range("F54").PasteSpecial(-4122)
$ws.Range("F58").PasteSpecial(-4122)

# Helper: write a cell as a literal text value without letting Excel's
# locale-aware auto-conversion turn ambiguous dd.mm.yyyy-looking strings
# (e.g. "10.10.2025") into date serials. Round-tripping through a quoted
# text formula and then pasting back as a value keeps the cell's existing
# number format / style untouched while still landing a plain text value.
function Set-TextValue($cell, $text) {
    $escaped = $text -replace '"', '""'
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

$rows = @(
    @{ Row = 54; A = "2026-02-17 17:05:28"; B = "poqui";           C = "Świeżo wykończone mieszkanie z dużym balkonem - Ponikwoda";                D = 2299;  E = "19.01.2026"; F = 29;  G = "https://www.olx.pl/d/oferta/swiezo-wykonczone-mieszkanie-z-duzym-balkonem-ponikwoda-CID3-ID1951OR.html";                H = "swiezo-wykonczone-mieszkanie-z-duzym-balkonem-ponikwoda-CID3-ID1951OR" }
    @{ Row = 55; A = "2026-02-17 17:05:28"; B = "poqui";           C = "Kawalerka po remoncie z funkcjonalną antresolą - ul. Jana Sawy";           D = 2499;  E = "28.10.2025"; F = 112; G = "https://www.olx.pl/d/oferta/kawalerka-po-remoncie-z-funkcjonalna-antresola-ul-jana-sawy-CID3-ID183ger.html";               H = "kawalerka-po-remoncie-z-funkcjonalna-antresola-ul-jana-sawy-CID3-ID183ger" }
    @{ Row = 56; A = "2026-02-17 17:05:28"; B = "poqui";           C = "Przytulny pokój blisko Politechniki – ul. Przytulna";                      D = 599;   E = "10.10.2025"; F = 130; G = "https://www.olx.pl/d/oferta/przytulny-pokoj-blisko-politechniki-ul-przytulna-CID3-ID17NeTz.html";                          H = "przytulny-pokoj-blisko-politechniki-ul-przytulna-CID3-ID17NeTz" }
    @{ Row = 57; A = "2026-02-17 17:05:28"; B = "pokojewlublinie"; C = "WOLNY OD ZARAZ! Pokój jedynka, ul. Romanowskiego 58";                      D = 58640; E = "11.08.2025"; F = 190; G = "https://www.olx.pl/d/oferta/wolny-od-zaraz-pokoj-jedynka-ul-romanowskiego-58-CID3-ID16ZeYm.html";                          H = "wolny-od-zaraz-pokoj-jedynka-ul-romanowskiego-58-CID3-ID16ZeYm" }
    @{ Row = 58; A = "2026-02-17 17:05:28"; B = "pokojewlublinie"; C = "WOLNY OD ZARAZ! Super lokalizacja, blisko centrum, ul. Paganiniego 12";    D = 12640; E = "19.01.2026"; F = 29;  G = "https://www.olx.pl/d/oferta/wolny-od-zaraz-super-lokalizacja-blisko-centrum-ul-paganiniego-12-CID3-ID195dLc.html";        H = "wolny-od-zaraz-super-lokalizacja-blisko-centrum-ul-paganiniego-12-CID3-ID195dLc" }
    @{ Row = 59; A = "2026-02-17 17:05:28"; B = "dawnypatron";     C = "Ładny pokój jednoosobowy. Wynajmę duży pokój w centrum. ul Niecała 4.";     D = 730;   E = "20.09.2024"; F = 515; G = "https://www.olx.pl/d/oferta/ladny-pokoj-jednoosobowy-wynajme-duzy-pokoj-w-centrum-ul-niecala-4-CID3-ID122jPM.html";         H = "ladny-pokoj-jednoosobowy-wynajme-duzy-pokoj-w-centrum-ul-niecala-4-CID3-ID122jPM" }
    @{ Row = 60; A = "2026-02-17 17:05:28"; B = "dawnypatron";     C = "Mam do wynajęcia pokój dla os. pracującej lub studenta. Narutowicza 14";    D = 14690; E = "05.12.2025"; F = 74;  G = "https://www.olx.pl/d/oferta/mam-do-wynajecia-pokoj-dla-os-pracujacej-lub-studenta-narutowicza-14-CID3-ID18ySfv.html";      H = "mam-do-wynajecia-pokoj-dla-os-pracujacej-lub-studenta-narutowicza-14-CID3-ID18ySfv" }
)

foreach ($r in $rows) {
    $i = $r.Row
    Set-TextValue $ws.Cells.Item($i, 1) $r.A
    Set-TextValue $ws.Cells.Item($i, 2) $r.B
    Set-TextValue $ws.Cells.Item($i, 3) $r.C
    $ws.Cells.Item($i, 4).Value = $r.D
    Set-TextValue $ws.Cells.Item($i, 5) $r.E
    $ws.Cells.Item($i, 6).Value = $r.F
    Set-TextValue $ws.Cells.Item($i, 7) $r.G
    Set-TextValue $ws.Cells.Item($i, 8) $r.H
}

Write-Output "Appended rows 54-60 to PODSUMOWANIE sheet"
